# Generate Report for Handback
# Updates the handback-status workbook so the "latest" handoff/handback
# file references (previously 30437b27-... / a95e265a-...) point at the
# newly generated files (d5c9d64f-... / ffff42586c1e-...), and refreshes
# the associated generation / handback timestamps.

$wb = $excel.ActiveWorkbook

$oldFile1 = "30437b27-e370-4827-8fba-ffc3720fe6b0"
$oldFile2 = "a95e265a-fe6b-4795-8d86-d0623c821d77"
$newFile1 = "d5c9d64f-e78d-4a3b-87b6-9fd4d40cafad"
$newFile2 = "ffff42586c1e-1d48-4b89-a744-20cbad1b1789"

$newFile1Md = "$newFile1.md"
$newFile2Md = "$newFile2.md"
$newFile1Path = "e2e\$newFile1.md"
$newFile2Path = "e2e\$newFile2.md"

$newXliffZh = "$newFile1.61cb262f4e41b915f460da185087fa1efba904c7.zh-cn.xlf"
$newXliffDe = "$newFile1.61cb262f4e41b915f460da185087fa1efba904c7.de-de.xlf"

$newHandoffDate = "2016-08-29 11:04:11"
$newZhHandoffDate = "2016-08-29 11:04:01"
$newZhHandbackDate = "2016-08-29 11:04:29"
$newDeHandbackDate = "2016-08-29 11:04:36"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1Md
$wsOverview.Range("B2").Value = $newFile1Path
$wsOverview.Range("G2").Value = $newHandoffDate

$wsOverview.Range("A3").Value = $newFile2Md
$wsOverview.Range("B3").Value = $newFile2Path
$wsOverview.Range("G3").Value = $newHandoffDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5949e19d16699cd6e1d67ec3cfa133c1c33ba222/e2e/$newFile1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile1Path
) | Out-Null
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5949e19d16699cd6e1d67ec3cfa133c1c33ba222/e2e/$newFile2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile2Path
) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFile1Md
$wsZh.Range("G2").Value = $newXliffZh
$wsZh.Range("H2").Value = $newZhHandoffDate
$wsZh.Range("I2").Value = $newFile1Md
$wsZh.Range("J2").Value = $newXliffZh
$wsZh.Range("K2").Value = $newZhHandbackDate

$wsZh.Range("A3").Value = $newFile2Md
$wsZh.Range("G3").Value = $newXliffZh
$wsZh.Range("H3").Value = $newZhHandoffDate
$wsZh.Range("I3").Value = $newFile2Md
$wsZh.Range("J3").Value = $newXliffZh
$wsZh.Range("K3").Value = $newZhHandbackDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5949e19d16699cd6e1d67ec3cfa133c1c33ba222/e2e/$newFile1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile1Md
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ed3a186413683ad8b6bd25f9f47b896802de3718/e2e/$newFile1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile1Md
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5949e19d16699cd6e1d67ec3cfa133c1c33ba222/e2e/$newFile2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile2Md
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ed3a186413683ad8b6bd25f9f47b896802de3718/e2e/$newFile2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile2Md
) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFile1Md
$wsDe.Range("G2").Value = $newXliffDe
$wsDe.Range("H2").Value = $newHandoffDate
$wsDe.Range("I2").Value = $newFile1Md
$wsDe.Range("J2").Value = $newXliffDe
$wsDe.Range("K2").Value = $newDeHandbackDate

$wsDe.Range("A3").Value = $newFile2Md
$wsDe.Range("G3").Value = $newXliffDe
$wsDe.Range("H3").Value = $newHandoffDate
$wsDe.Range("I3").Value = $newFile2Md
$wsDe.Range("J3").Value = $newXliffDe
$wsDe.Range("K3").Value = $newDeHandbackDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5949e19d16699cd6e1d67ec3cfa133c1c33ba222/e2e/$newFile1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile1Md
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8ee768fcda45417cb467d5e19800bdc2eef3e82f/e2e/$newFile1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile1Md
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5949e19d16699cd6e1d67ec3cfa133c1c33ba222/e2e/$newFile2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile2Md
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8ee768fcda45417cb467d5e19800bdc2eef3e82f/e2e/$newFile2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile2Md
) | Out-Null

Write-Host "Handback status report regenerated."
